$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date
$ws.Name = "Through 2022-06-24"

# Update the header cell text (shared string) that also carries the "through" date
$ws.Range("I1").Value = "2022 (through 06-24)"

# Update the new data values for July (row 7) and Total (row 14)
$ws.Range("I7").Value = 114
$ws.Range("I14").Value = 777
